$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$s.Shapes.Item(25).Delete()
